# Reorder the "secciones" columns: insert a brand-new question column
# before column I ("preg_secc0") and push the existing I:L question
# columns (and their header comments) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the current comment texts before anything shifts.
$txt_I = $ws.Range("I1").Comment.Text()
$txt_J = $ws.Range("J1").Comment.Text()
$txt_K = $ws.Range("K1").Comment.Text()
$txt_L = $ws.Range("L1").Comment.Text()

# Insert a new column before I; this shifts columns I..L -> J..M
# (cell values, including the sparse I/J/K data columns, move with it).
$ws.Range("I1").EntireColumn.Insert()

# New header for the freshly inserted column.
$ws.Range("I1").Value = "preg_secc0"

# Re-home the old comments on the cells they now occupy (comments do not
# follow the column insert automatically), preserving original order.
$ws.Range("M1").AddComment($txt_L)
[void]$ws.Range("L1").Comment.Text($txt_K)
[void]$ws.Range("K1").Comment.Text($txt_J)
[void]$ws.Range("J1").Comment.Text($txt_I)
$ws.Range("I1").AddComment("testing secciones (seccion 0, debería ir primero)")
